$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add 3 new groups of data (cxq in E, hyy in G, hzj in F) with their headers
$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("G1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("F1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

$ws.Range("E2").Value = 0.99109792284866471
$ws.Range("F2").Value = 0.963963963963964
$ws.Range("G2").Value = 0.9494949494949495

$ws.Range("E3").Value = 0.95918367346938771
$ws.Range("F3").Value = 0.93548387096774199
$ws.Range("G3").Value = 0.94880546075085326

# Update the current selection to match the committed view state
$ws.Range("F1:F1048576").Select()

# Resize the workbook window to match the saved view state
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 13170
